# WorkReport.xlsx edit: add four new work-log entries (rows 65-68) to
# Sheet1, covering FillFormReport/CSV+GoogleDocs export, the Wiki .NET
# parser research, MarkItUp integration, and script editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the description column in the order the original author appears to
# have typed them (row65, row67, row66, row68) so the workbook's shared-
# string table picks up the four new unique strings in that same order:
#   60 FillFormReport, Export do CSV a Google Docs
#   61 Integrace MarkItUp
#   62 Zjistovani moznosti Wiki (.net, JS), Skripta - DB, modely atd, Integrace Wiki .NET parser
#   63 Editace skript
$ws.Range("B65").Value = "FillFormReport, Export do CSV a Google Docs"
$ws.Range("B67").Value = "Integrace MarkItUp"
$ws.Range("B66").Value = "Zjistovani moznosti Wiki (.net, JS), Skripta - DB, modely atd, Integrace Wiki .NET parser"
$ws.Range("B68").Value = "Editace skript"

# Hours worked (column C) and the work date as an OLE Automation date
# serial (column D) for each new row.
$ws.Range("C65").Value = 8
$ws.Range("D65").Value = 40999   # 2012-03-31

$ws.Range("C66").Value = 8
$ws.Range("D66").Value = 41000   # 2012-04-01

$ws.Range("C67").Value = 2
$ws.Range("D67").Value = 41001   # 2012-04-02

$ws.Range("C68").Value = 2
$ws.Range("D68").Value = 41002   # 2012-04-03

# Give the new date cells the same date-formatted style already used by the
# rest of column D (copy/paste-special formats reuses the existing style
# instead of minting a new numFmt).
$ws.Range("D64").Copy()
$ws.Range("D65:D68").PasteSpecial(-4122)

# Land the selection on the next empty row, same as Excel leaves the cursor
# after entering the last row of data.
$ws.Range("B69").Select()
$ws.Application.ActiveWindow.ScrollRow = 35
